$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text (string) representation
# instead of being auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '30.401.04'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '1.916.44'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '241.08'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = '0.4671'
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('D8').Value = '0.2841'
$ws.Range('D9').Value = '0.06823'
$ws.Range('E9').Value = '  +4.55%  '
$ws.Range('D10').Value = '106.92'
$ws.Range('E10').Value = '  +12.90%  '
$ws.Range('D11').Value = '17.94'
$ws.Range('E11').Value = '  -5.35%  '
$ws.Range('D12').Value = '1.912.07'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '0.07617'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').Value = '0.6521'
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').Value = '286.14'
$ws.Range('E16').Value = '  -4.29%  '
$ws.Range('D17').Value = '30.408.41'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').Value = '0.000007571'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').Value = '2.161.13'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = '5.192'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').Value = '168.07'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '9.223'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').Value = '21.33'
$ws.Range('E27').Value = '  +8.49%  '
$ws.Range('D28').Value = '2.028'
$ws.Range('E28').Value = '  +3.80%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = '4.125'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').Value = '3.934'
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').Value = '0.05022'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = '0.7350'
$ws.Range('E34').Value = '  +1.32%  '
$ws.Range('D35').Value = '1.141'
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('D36').Value = '0.9991'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '2.728'
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').Value = '0.02008'
$ws.Range('E38').Value = '  +2.93%  '
$ws.Range('D39').Value = '2.678'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = '2.038'
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('D41').Value = '108.73'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').Value = '0.8714'
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('D43').Value = '5.804'
$ws.Range('E43').Value = '  +3.61%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '53.01'
$ws.Range('E44').Value = '  +25.91%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9995'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('D48').Value = '7.109'
$ws.Range('E48').Value = '  -3.54%  '
$ws.Range('D49').Value = '9.104'
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('D51').Value = '34.47'
$ws.Range('E51').Value = '  -0.57%  '
